$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70: extend text with two more lines ("Kurzes Teammeeting" / "Statusupdate
# erstellt"), bump hours 2 -> 5. Do this before touching row 66 so the extended
# string is registered (in the shared-strings table) ahead of the TdOT text below.
$ws.Cells.Item(70, 7).Value = "Optimierung der stream UI für die neu implementierten Sprachen Französisch und Spanisch`nErgänzungen der ResourceBundles`nDokumentation der noch ausstehenden Funktionen`nKurzes Teammeeting`nStatusupdate erstellt"
$ws.Cells.Item(70, 6).Value = 5
$ws.Rows.Item(70).RowHeight = 72

# Row 66: mention the new 'Über' page explicitly in the TdOT note
$ws.Cells.Item(66, 7).Value = "Codedokumentierung`nAnpassungen für TdOT (Einfügen neuer Seite 'Über', Counter bis Release)`nDiverse Vorbereitungen für TdOT getroffen"

# Row 71: brand-new entry for the css-files work.
$ws.Cells.Item(71, 5).Value = 43844
$ws.Cells.Item(71, 6).Value = 5.5
$ws.Cells.Item(71, 7).Value = "Arbeiten an stream App Icon`nEinführung von css files (Notwendig durch die zukünftige Implementierung eines Dark modes)`nAnnpassungen an stream UI wegen Einführung von css files"

# Copy row 70's formatting (date / number / wrap-text styles) onto the new row so
# it matches the existing look; values were already set above so the SUM(F:F)
# formula in C5 stays correctly wired to the new F71 cell.
$ws.Range("E70:G70").Copy()
$ws.Range("E71:G71").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(71).RowHeight = 43.2

$ws.Range("F69").Select() | Out-Null
